{"js": "// Update the date line and all 25 \"two-digit \u00d7 two-digit\" answer cells\n// in the practice-sheet table to the values from the next day's sheet.\nconst replacements = [\n  [\"2024-11-13 Wednesday\", \"2024-11-14 Thursday\"],\n  [\"44\u00d783=3652\", \"83\u00d791=7553\"],\n  [\"37\u00d791=3367\", \"57\u00d727=1539\"],\n  [\"55\u00d771=3905\", \"61\u00d778=4758\"],\n  [\"42\u00d758=2436\", \"33\u00d795=3135\"],\n  [\"88\u00d777=6776\", \"50\u00d721=1050\"],\n  [\"77\u00d794=7238\", \"87\u00d714=1218\"],\n  [\"11\u00d788=968\", \"90\u00d752=4680\"],\n  [\"96\u00d759=5664\", \"27\u00d732=864\"],\n  [\"33\u00d768=2244\", \"71\u00d737=2627\"],\n  [\"95\u00d791=8645\", \"29\u00d769=2001\"],\n  [\"90\u00d759=5310\", \"17\u00d762=1054\"],\n  [\"26\u00d798=2548\", \"63\u00d794=5922\"],\n  [\"37\u00d787=3219\", \"31\u00d786=2666\"],\n  [\"25\u00d769=1725\", \"78\u00d740=3120\"],\n  [\"31\u00d730=930\", \"50\u00d723=1150\"],\n  [\"15\u00d712=180\", \"28\u00d740=1120\"],\n  [\"19\u00d799=1881\", \"19\u00d744=836\"],\n  [\"21\u00d746=966\", \"99\u00d719=1881\"],\n  [\"47\u00d784=3948\", \"12\u00d780=960\"],\n  [\"39\u00d742=1638\", \"92\u00d753=4876\"],\n  [\"20\u00d753=1060\", \"85\u00d742=3570\"],\n  [\"53\u00d775=3975\", \"29\u00d741=1189\"],\n  [\"93\u00d728=2604\", \"87\u00d742=3654\"],\n  [\"26\u00d765=1690\", \"44\u00d777=3388\"],\n  [\"76\u00d745=3420\", \"65\u00d795=6175\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and all 25 \"two-digit \u00d7 two-digit\" answer cells\n# in the practice-sheet table to the values from the next day's sheet.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-11-13 Wednesday\", \"2024-11-14 Thursday\"),\n  @(\"44\u00d783=3652\", \"83\u00d791=7553\"),\n  @(\"37\u00d791=3367\", \"57\u00d727=1539\"),\n  @(\"55\u00d771=3905\", \"61\u00d778=4758\"),\n  @(\"42\u00d758=2436\", \"33\u00d795=3135\"),\n  @(\"88\u00d777=6776\", \"50\u00d721=1050\"),\n  @(\"77\u00d794=7238\", \"87\u00d714=1218\"),\n  @(\"11\u00d788=968\", \"90\u00d752=4680\"),\n  @(\"96\u00d759=5664\", \"27\u00d732=864\"),\n  @(\"33\u00d768=2244\", \"71\u00d737=2627\"),\n  @(\"95\u00d791=8645\", \"29\u00d769=2001\"),\n  @(\"90\u00d759=5310\", \"17\u00d762=1054\"),\n  @(\"26\u00d798=2548\", \"63\u00d794=5922\"),\n  @(\"37\u00d787=3219\", \"31\u00d786=2666\"),\n  @(\"25\u00d769=1725\", \"78\u00d740=3120\"),\n  @(\"31\u00d730=930\", \"50\u00d723=1150\"),\n  @(\"15\u00d712=180\", \"28\u00d740=1120\"),\n  @(\"19\u00d799=1881\", \"19\u00d744=836\"),\n  @(\"21\u00d746=966\", \"99\u00d719=1881\"),\n  @(\"47\u00d784=3948\", \"12\u00d780=960\"),\n  @(\"39\u00d742=1638\", \"92\u00d753=4876\"),\n  @(\"20\u00d753=1060\", \"85\u00d742=3570\"),\n  @(\"53\u00d775=3975\", \"29\u00d741=1189\"),\n  @(\"93\u00d728=2604\", \"87\u00d742=3654\"),\n  @(\"26\u00d765=1690\", \"44\u00d777=3388\"),\n  @(\"76\u00d745=3420\", \"65\u00d795=6175\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
